# Raw and Clean Data from SSA for June 26th (2020-06-26, Excel serial 44008)
#
# Appends the next day's row of data (44008) to the four daily-tracking
# sheets (out_vars, dates_dx, dates_sx, dates_deaths), fills in the new
# "AA" column for 2020-06-26 on control_obs (extending the running SUM
# formula along with it), and leaves the final cell selection on each
# sheet matching where the original author ended up while doing the data
# entry.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------
# 1) out_vars - brand new row 27 (A:J), formatted like row 26 above it.
# ---------------------------------------------------------------------
$wsOut = $wb.Worksheets.Item("out_vars")

$wsOut.Range("A26:J26").Copy() | Out-Null
$wsOut.Range("A27").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$wsOut.Range("A27").Value = 44008
$wsOut.Range("B27").Value = 208392
$wsOut.Range("C27").Value = 267288
$wsOut.Range("D27").Value = 66440
$wsOut.Range("E27").Value = 25779
$wsOut.Range("F27").Value = 31.374524933778648
$wsOut.Range("G27").Value = 65382
$wsOut.Range("H27").Value = 5570
$wsOut.Range("I27").Value = 5844
$wsOut.Range("J27").Value = 542120

# ---------------------------------------------------------------------
# 2) dates_dx - row 27 already exists (blank placeholders); fill it in.
# ---------------------------------------------------------------------
$wsDx = $wb.Worksheets.Item("dates_dx")

$wsDx.Range("A26").Copy() | Out-Null
$wsDx.Range("A27").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$wsDx.Range("A27").Value = 44008
$wsDx.Range("B27").Value = 0
$wsDx.Range("C27").Value = 1
$wsDx.Range("D27").Value = 1
$wsDx.Range("E27").Value = 1
$wsDx.Range("F27").Value = 1
$wsDx.Range("G27").Value = 0
$wsDx.Range("H27").Value = 0
$wsDx.Range("I27").Value = 1
$wsDx.Range("J27").Value = 0
$wsDx.Range("K27").Value = 4

# ---------------------------------------------------------------------
# 3) dates_sx - row 27 already exists (blank placeholders); fill it in.
# ---------------------------------------------------------------------
$wsSx = $wb.Worksheets.Item("dates_sx")

$wsSx.Range("A26").Copy() | Out-Null
$wsSx.Range("A27").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$wsSx.Range("A27").Value = 44008
$wsSx.Range("B27").Value = 0
$wsSx.Range("C27").Value = 1
$wsSx.Range("D27").Value = 1
$wsSx.Range("E27").Value = 0
$wsSx.Range("F27").Value = 1
$wsSx.Range("G27").Value = 1
$wsSx.Range("H27").Value = 1
$wsSx.Range("I27").Value = 0
$wsSx.Range("J27").Value = 1
$wsSx.Range("K27").Value = 1
$wsSx.Range("L27").Value = 0
$wsSx.Range("M27").Value = 0

# ---------------------------------------------------------------------
# 4) dates_deaths - row 27 already exists (blank placeholders); fill it in.
# ---------------------------------------------------------------------
$wsDeaths = $wb.Worksheets.Item("dates_deaths")

$wsDeaths.Range("A26").Copy() | Out-Null
$wsDeaths.Range("A27").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$wsDeaths.Range("A27").Value = 44008
$wsDeaths.Range("B27").Value = 0
$wsDeaths.Range("C27").Value = 0
$wsDeaths.Range("D27").Value = 0
$wsDeaths.Range("E27").Value = 0
$wsDeaths.Range("F27").Value = 2
$wsDeaths.Range("G27").Value = 1
$wsDeaths.Range("H27").Value = 1
$wsDeaths.Range("I27").Value = 1
$wsDeaths.Range("J27").Value = 2

# ---------------------------------------------------------------------
# 5) control_obs - new column AA for 2020-06-26, plus extending the
#    running-total SUM formula in row 20 one column to the right.
# ---------------------------------------------------------------------
$wsControl = $wb.Worksheets.Item("control_obs")

$wsControl.Range("Z1").Copy() | Out-Null
$wsControl.Range("AA1").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$wsControl.Range("AA1").Value = 44008
$wsControl.Range("AA2").Value = 3849
$wsControl.Range("AA3").Value = 3664
$wsControl.Range("AA4").Value = 3664
$wsControl.Range("AA5").Value = 3664
$wsControl.Range("AA6").Value = 3664
$wsControl.Range("AA7").Value = 2855
$wsControl.Range("AA8").Value = 5486
$wsControl.Range("AA10").Value = 166
$wsControl.Range("AA11").Value = 166
$wsControl.Range("AA12").Value = 166
$wsControl.Range("AA13").Value = 166
$wsControl.Range("AA14").Value = 166
$wsControl.Range("AA15").Value = 101
$wsControl.Range("AA16").Value = 178
$wsControl.Range("AA18").Value = 891

$wsControl.Range("Z20").Copy() | Out-Null
$wsControl.Range("AA20").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false
$wsControl.Range("AA20").Formula = "=SUM(AA2:AA18)"

# ---------------------------------------------------------------------
# 6) Selections - restore the cell each sheet was left on, matching the
#    author's final view. Finish on out_vars so it stays the active tab.
# ---------------------------------------------------------------------
$wsDx.Activate() | Out-Null
$wsDx.Range("C31").Select() | Out-Null

$wsSx.Activate() | Out-Null
$wsSx.Range("A28").Select() | Out-Null

$wsDeaths.Activate() | Out-Null
$wsDeaths.Range("E24").Select() | Out-Null

$wsControlMpio = $wb.Worksheets.Item("control_obs_mpio")
$wsControlMpio.Activate() | Out-Null
$wsControlMpio.Range("G8").Select() | Out-Null

$wsControl.Activate() | Out-Null
$wsControl.Range("X19").Select() | Out-Null

$wsAnomalias = $wb.Worksheets.Item("anomalias")
$wsAnomalias.Activate() | Out-Null
$wsAnomalias.Range("D12").Select() | Out-Null

$wsOut.Activate() | Out-Null
$wsOut.Range("D12").Select() | Out-Null
